$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a brand-new first paragraph containing "Change1" and move
#    the "_GoBack" bookmark onto it (right after the new text, still
#    inside the same paragraph).
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs.Item(1).Range.Text = "Change1"

# Bookmark covering exactly the "Change1" text (not the paragraph mark).
# Adding a bookmark with the same name as an existing one moves it here
# (Word enforces unique bookmark names), so the old "_GoBack" bookmark
# that sat in its own empty paragraph further down is removed from
# there automatically.
$newBookmarkRange = $d.Range(0, 7)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null

# ------------------------------------------------------------------
# 2) In the "Infoway (...)" paragraph, drop the grammar-check markers
#    and merge the trailing "Implementing partner)->", "Cisco(",
#    "end client)." runs into a single run, while leaving the
#    preceding "America (" run untouched.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Infoway (*end client).*") {
        $infowayPara = $candidate
        break
    }
}

$paraText = $infowayPara.Range.Text
$target = "Implementing partner)->Cisco(end client)."
$offset = $paraText.IndexOf($target)
$mergeStart = $infowayPara.Range.Start + $offset
$mergeEnd = $mergeStart + $target.Length
$mergeRange = $d.Range($mergeStart, $mergeEnd)

$mergeXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Implementing partner)-&gt;Cisco(end client).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$mergeRange.InsertXML($mergeXml)
